$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain numeric updates (values that genuinely changed, no sign-of-zero subtlety)
$ws.Range("L2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("J14").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("J17").Value = 2
$ws.Range("L17").Value = 2
$ws.Range("M17").Value = 1
$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 11

# Cells whose solver output is a floating negative zero (-0.0): same magnitude as 0,
# but we reproduce the exact IEEE -0.0 bit pattern via a scratch formula + values-only paste,
# since a direct Value assignment of a literal negative zero normalizes to +0.
$ws.Range("ZZ1").Formula = "=-1*0"
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").PasteSpecial(-4163)
$ws.Range("F2").PasteSpecial(-4163)
$ws.Range("J2").PasteSpecial(-4163)
$ws.Range("M2").PasteSpecial(-4163)
$ws.Range("P2").PasteSpecial(-4163)
$ws.Range("C3").PasteSpecial(-4163)
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").PasteSpecial(-4163)
$ws.Range("F3").PasteSpecial(-4163)
$ws.Range("H3").PasteSpecial(-4163)
$ws.Range("J3").PasteSpecial(-4163)
$ws.Range("M3").PasteSpecial(-4163)
$ws.Range("P3").PasteSpecial(-4163)
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("F4").PasteSpecial(-4163)
$ws.Range("H4").PasteSpecial(-4163)
$ws.Range("J4").PasteSpecial(-4163)
$ws.Range("L4").PasteSpecial(-4163)
$ws.Range("P4").PasteSpecial(-4163)
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("C5").PasteSpecial(-4163)
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").PasteSpecial(-4163)
$ws.Range("F5").PasteSpecial(-4163)
$ws.Range("H5").PasteSpecial(-4163)
$ws.Range("L5").PasteSpecial(-4163)
$ws.Range("M5").PasteSpecial(-4163)
$ws.Range("B6").PasteSpecial(-4163)
$ws.Range("C6").PasteSpecial(-4163)
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").PasteSpecial(-4163)
$ws.Range("H6").PasteSpecial(-4163)
$ws.Range("J6").PasteSpecial(-4163)
$ws.Range("L6").PasteSpecial(-4163)
$ws.Range("M6").PasteSpecial(-4163)
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("F8").PasteSpecial(-4163)
$ws.Range("H8").PasteSpecial(-4163)
$ws.Range("L8").PasteSpecial(-4163)
$ws.Range("B9").PasteSpecial(-4163)
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("B10").PasteSpecial(-4163)
$ws.Range("C10").PasteSpecial(-4163)
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("F10").PasteSpecial(-4163)
$ws.Range("L10").PasteSpecial(-4163)
$ws.Range("B11").PasteSpecial(-4163)
$ws.Range("C11").PasteSpecial(-4163)
$ws.Range("H11").PasteSpecial(-4163)
$ws.Range("L11").PasteSpecial(-4163)
$ws.Range("P12").PasteSpecial(-4163)
$ws.Range("M14").PasteSpecial(-4163)
$ws.Range("P14").PasteSpecial(-4163)
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("M15").PasteSpecial(-4163)
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("J16").PasteSpecial(-4163)
$ws.Range("M16").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()
$excel.CutCopyMode = $false
